$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-08 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-09 Friday", 2) | Out-Null
$d.Content.Find.Execute("302÷5=60, 2", $true, $false, $false, $false, $false, $true, 1, $false, "730÷6=121, 4", 2) | Out-Null
$d.Content.Find.Execute("207÷9=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "890÷6=148, 2", 2) | Out-Null
$d.Content.Find.Execute("864÷5=172, 4", $true, $false, $false, $false, $false, $true, 1, $false, "690÷6=115, 0", 2) | Out-Null
$d.Content.Find.Execute("144÷2=72, 0", $true, $false, $false, $false, $false, $true, 1, $false, "270÷3=90, 0", 2) | Out-Null
$d.Content.Find.Execute("640÷2=320, 0", $true, $false, $false, $false, $false, $true, 1, $false, "513÷9=57, 0", 2) | Out-Null
$d.Content.Find.Execute("303÷8=37, 7", $true, $false, $false, $false, $false, $true, 1, $false, "133÷9=14, 7", 2) | Out-Null
$d.Content.Find.Execute("863÷7=123, 2", $true, $false, $false, $false, $false, $true, 1, $false, "693÷9=77, 0", 2) | Out-Null
$d.Content.Find.Execute("397÷4=99, 1", $true, $false, $false, $false, $false, $true, 1, $false, "495÷5=99, 0", 2) | Out-Null
$d.Content.Find.Execute("455÷6=75, 5", $true, $false, $false, $false, $false, $true, 1, $false, "572÷2=286, 0", 2) | Out-Null
$d.Content.Find.Execute("743÷8=92, 7", $true, $false, $false, $false, $false, $true, 1, $false, "149÷6=24, 5", 2) | Out-Null
$d.Content.Find.Execute("683÷7=97, 4", $true, $false, $false, $false, $false, $true, 1, $false, "132÷8=16, 4", 2) | Out-Null
$d.Content.Find.Execute("106÷7=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "985÷3=328, 1", 2) | Out-Null
$d.Content.Find.Execute("940÷7=134, 2", $true, $false, $false, $false, $false, $true, 1, $false, "389÷7=55, 4", 2) | Out-Null
$d.Content.Find.Execute("957÷5=191, 2", $true, $false, $false, $false, $false, $true, 1, $false, "169÷8=21, 1", 2) | Out-Null
$d.Content.Find.Execute("799÷2=399, 1", $true, $false, $false, $false, $false, $true, 1, $false, "755÷2=377, 1", 2) | Out-Null
$d.Content.Find.Execute("266÷9=29, 5", $true, $false, $false, $false, $false, $true, 1, $false, "198÷4=49, 2", 2) | Out-Null
$d.Content.Find.Execute("994÷2=497, 0", $true, $false, $false, $false, $false, $true, 1, $false, "155÷2=77, 1", 2) | Out-Null
$d.Content.Find.Execute("354÷5=70, 4", $true, $false, $false, $false, $false, $true, 1, $false, "959÷3=319, 2", 2) | Out-Null
$d.Content.Find.Execute("525÷6=87, 3", $true, $false, $false, $false, $false, $true, 1, $false, "182÷5=36, 2", 2) | Out-Null
$d.Content.Find.Execute("587÷3=195, 2", $true, $false, $false, $false, $false, $true, 1, $false, "900÷7=128, 4", 2) | Out-Null
$d.Content.Find.Execute("414÷9=46, 0", $true, $false, $false, $false, $false, $true, 1, $false, "321÷2=160, 1", 2) | Out-Null
$d.Content.Find.Execute("642÷4=160, 2", $true, $false, $false, $false, $false, $true, 1, $false, "178÷6=29, 4", 2) | Out-Null
$d.Content.Find.Execute("440÷4=110, 0", $true, $false, $false, $false, $false, $true, 1, $false, "609÷5=121, 4", 2) | Out-Null
$d.Content.Find.Execute("694÷8=86, 6", $true, $false, $false, $false, $false, $true, 1, $false, "923÷5=184, 3", 2) | Out-Null
$d.Content.Find.Execute("810÷3=270, 0", $true, $false, $false, $false, $false, $true, 1, $false, "683÷8=85, 3", 2) | Out-Null
